# Update NATMI LR-pair TPM data for Slit1-Robo1: add "ECs" as a 5th
# sending/target cluster, re-sort cluster rows alphabetically, and
# refresh all computed expression/specificity values (rows 2-21).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Slit1"
$ws.Cells.Item(2,3).Value = "Robo1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.08454100000000002
$ws.Cells.Item(2,8).Value = 0.253623
$ws.Cells.Item(2,9).Value = 0.04188307112135965
$ws.Cells.Item(2,10).Value = 0.04356611262330446
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.169654
$ws.Cells.Item(2,14).Value = 0.508962
$ws.Cells.Item(2,15).Value = 0.006094264463659866
$ws.Cells.Item(2,16).Value = 0.006534681579452628
$ws.Cells.Item(2,17).Value = 0.014342718814
$ws.Cells.Item(2,18).Value = 0.129084469326
$ws.Cells.Item(2,19).Value = 0.0002552465119638409
$ws.Cells.Item(2,20).Value = 0.0002846906736478663

# Row 3: ECs -> FAPs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Slit1"
$ws.Cells.Item(3,3).Value = "Robo1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.08454100000000002
$ws.Cells.Item(3,8).Value = 0.253623
$ws.Cells.Item(3,9).Value = 0.04188307112135965
$ws.Cells.Item(3,10).Value = 0.04356611262330446
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 21.99231
$ws.Cells.Item(3,14).Value = 65.97693
$ws.Cells.Item(3,15).Value = 0.7900017288527916
$ws.Cells.Item(3,16).Value = 0.8470931604713817
$ws.Cells.Item(3,17).Value = 1.85925187971
$ws.Cells.Item(3,18).Value = 16.73326691739
$ws.Cells.Item(3,19).Value = 0.03308769859553855
$ws.Cells.Item(3,20).Value = 0.03690455603152713

# Row 4: ECs -> Inflammatory-Mac
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Slit1"
$ws.Cells.Item(4,3).Value = "Robo1"
$ws.Cells.Item(4,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.08454100000000002
$ws.Cells.Item(4,8).Value = 0.253623
$ws.Cells.Item(4,9).Value = 0.04188307112135965
$ws.Cells.Item(4,10).Value = 0.04356611262330446
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.04769766666666667
$ws.Cells.Item(4,14).Value = 0.143093
$ws.Cells.Item(4,15).Value = 0.001713382501834088
$ws.Cells.Item(4,16).Value = 0.001837204332049573
$ws.Cells.Item(4,17).Value = 0.004032408437666667
$ws.Cells.Item(4,18).Value = 0.03629167593900001
$ws.Cells.Item(4,19).Value = 0.00007176172118241025
$ws.Cells.Item(4,20).Value = 0.00008003985084209454

# Row 5: ECs -> MuSCs
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Slit1"
$ws.Cells.Item(5,3).Value = "Robo1"
$ws.Cells.Item(5,4).Value = "MuSCs"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.08454100000000002
$ws.Cells.Item(5,8).Value = 0.253623
$ws.Cells.Item(5,9).Value = 0.04188307112135965
$ws.Cells.Item(5,10).Value = 0.04356611262330446
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 5.6286445
$ws.Cells.Item(5,14).Value = 11.257289
$ws.Cells.Item(5,15).Value = 0.2021906241817143
$ws.Cells.Item(5,16).Value = 0.1445349536171162
$ws.Cells.Item(5,17).Value = 0.4758512346745001
$ws.Cells.Item(5,18).Value = 2.855107408047
$ws.Cells.Item(5,19).Value = 0.00846836429267484
$ws.Cells.Item(5,20).Value = 0.00629682606728737

# Row 6: FAPs -> ECs
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Slit1"
$ws.Cells.Item(6,3).Value = "Robo1"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.340305
$ws.Cells.Item(6,8).Value = 1.020915
$ws.Cells.Item(6,9).Value = 0.1685929728528678
$ws.Cells.Item(6,10).Value = 0.1753677618702597
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.169654
$ws.Cells.Item(6,14).Value = 0.508962
$ws.Cells.Item(6,15).Value = 0.006094264463659866
$ws.Cells.Item(6,16).Value = 0.006534681579452628
$ws.Cells.Item(6,17).Value = 0.05773410447000001
$ws.Cells.Item(6,18).Value = 0.51960694023
$ws.Cells.Item(6,19).Value = 0.001027450163280004
$ws.Cells.Item(6,20).Value = 0.001145972483123421

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Slit1"
$ws.Cells.Item(7,3).Value = "Robo1"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.340305
$ws.Cells.Item(7,8).Value = 1.020915
$ws.Cells.Item(7,9).Value = 0.1685929728528678
$ws.Cells.Item(7,10).Value = 0.1753677618702597
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 21.99231
$ws.Cells.Item(7,14).Value = 65.97693
$ws.Cells.Item(7,15).Value = 0.7900017288527916
$ws.Cells.Item(7,16).Value = 0.8470931604713817
$ws.Cells.Item(7,17).Value = 7.484093054550001
$ws.Cells.Item(7,18).Value = 67.35683749095
$ws.Cells.Item(7,19).Value = 0.1331887400261973
$ws.Cells.Item(7,20).Value = 0.1485528316474709

# Row 8: FAPs -> Inflammatory-Mac
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Slit1"
$ws.Cells.Item(8,3).Value = "Robo1"
$ws.Cells.Item(8,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.340305
$ws.Cells.Item(8,8).Value = 1.020915
$ws.Cells.Item(8,9).Value = 0.1685929728528678
$ws.Cells.Item(8,10).Value = 0.1753677618702597
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.04769766666666667
$ws.Cells.Item(8,14).Value = 0.143093
$ws.Cells.Item(8,15).Value = 0.001713382501834088
$ws.Cells.Item(8,16).Value = 0.001837204332049573
$ws.Cells.Item(8,17).Value = 0.016231754455
$ws.Cells.Item(8,18).Value = 0.146085790095
$ws.Cells.Item(8,19).Value = 0.0002888642496182931
$ws.Cells.Item(8,20).Value = 0.000322186411809879

# Row 9: FAPs -> MuSCs
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Slit1"
$ws.Cells.Item(9,3).Value = "Robo1"
$ws.Cells.Item(9,4).Value = "MuSCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.340305
$ws.Cells.Item(9,8).Value = 1.020915
$ws.Cells.Item(9,9).Value = 0.1685929728528678
$ws.Cells.Item(9,10).Value = 0.1753677618702597
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 5.6286445
$ws.Cells.Item(9,14).Value = 11.257289
$ws.Cells.Item(9,15).Value = 0.2021906241817143
$ws.Cells.Item(9,16).Value = 0.1445349536171162
$ws.Cells.Item(9,17).Value = 1.9154558665725
$ws.Cells.Item(9,18).Value = 11.492735199435
$ws.Cells.Item(9,19).Value = 0.03408791841377214
$ws.Cells.Item(9,20).Value = 0.02534677132785546

# Row 10: Inflammatory-Mac -> ECs
$ws.Cells.Item(10,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10,2).Value = "Slit1"
$ws.Cells.Item(10,3).Value = "Robo1"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.4435069999999999
$ws.Cells.Item(10,8).Value = 1.330521
$ws.Cells.Item(10,9).Value = 0.2197210255830999
$ws.Cells.Item(10,10).Value = 0.2285503591301722
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.169654
$ws.Cells.Item(10,14).Value = 0.508962
$ws.Cells.Item(10,15).Value = 0.006094264463659866
$ws.Cells.Item(10,16).Value = 0.006534681579452628
$ws.Cells.Item(10,17).Value = 0.07524273657799999
$ws.Cells.Item(10,18).Value = 0.677184629202
$ws.Cells.Item(10,19).Value = 0.001339038038129986
$ws.Cells.Item(10,20).Value = 0.001493503821785219

# Row 11: Inflammatory-Mac -> FAPs
$ws.Cells.Item(11,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(11,2).Value = "Slit1"
$ws.Cells.Item(11,3).Value = "Robo1"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 0.4435069999999999
$ws.Cells.Item(11,8).Value = 1.330521
$ws.Cells.Item(11,9).Value = 0.2197210255830999
$ws.Cells.Item(11,10).Value = 0.2285503591301722
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 21.99231
$ws.Cells.Item(11,14).Value = 65.97693
$ws.Cells.Item(11,15).Value = 0.7900017288527916
$ws.Cells.Item(11,16).Value = 0.8470931604713817
$ws.Cells.Item(11,17).Value = 9.753743431169998
$ws.Cells.Item(11,18).Value = 87.78369088052999
$ws.Cells.Item(11,19).Value = 0.1735799900759574
$ws.Cells.Item(11,20).Value = 0.1936034460424468

# Row 12: Inflammatory-Mac -> Inflammatory-Mac
$ws.Cells.Item(12,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(12,2).Value = "Slit1"
$ws.Cells.Item(12,3).Value = "Robo1"
$ws.Cells.Item(12,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 0.4435069999999999
$ws.Cells.Item(12,8).Value = 1.330521
$ws.Cells.Item(12,9).Value = 0.2197210255830999
$ws.Cells.Item(12,10).Value = 0.2285503591301722
$ws.Cells.Item(12,11).Value = 1
$ws.Cells.Item(12,12).Value = 0.3333333333333333
$ws.Cells.Item(12,13).Value = 0.04769766666666667
$ws.Cells.Item(12,14).Value = 0.143093
$ws.Cells.Item(12,15).Value = 0.001713382501834088
$ws.Cells.Item(12,16).Value = 0.001837204332049573
$ws.Cells.Item(12,17).Value = 0.02115424905033333
$ws.Cells.Item(12,18).Value = 0.190388241453
$ws.Cells.Item(12,19).Value = 0.0003764661605191234
$ws.Cells.Item(12,20).Value = 0.000419893709885438

# Row 13: Inflammatory-Mac -> MuSCs
$ws.Cells.Item(13,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(13,2).Value = "Slit1"
$ws.Cells.Item(13,3).Value = "Robo1"
$ws.Cells.Item(13,4).Value = "MuSCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 0.4435069999999999
$ws.Cells.Item(13,8).Value = 1.330521
$ws.Cells.Item(13,9).Value = 0.2197210255830999
$ws.Cells.Item(13,10).Value = 0.2285503591301722
$ws.Cells.Item(13,11).Value = 2
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 5.6286445
$ws.Cells.Item(13,14).Value = 11.257289
$ws.Cells.Item(13,15).Value = 0.2021906241817143
$ws.Cells.Item(13,16).Value = 0.1445349536171162
$ws.Cells.Item(13,17).Value = 2.4963432362615
$ws.Cells.Item(13,18).Value = 14.978059417569
$ws.Cells.Item(13,19).Value = 0.04442553130849338
$ws.Cells.Item(13,20).Value = 0.03303351555605468

# Row 14: MuSCs -> ECs
$ws.Cells.Item(14,1).Value = "MuSCs"
$ws.Cells.Item(14,2).Value = "Slit1"
$ws.Cells.Item(14,3).Value = "Robo1"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 2
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 0.2339355
$ws.Cells.Item(14,8).Value = 0.467871
$ws.Cells.Item(14,9).Value = 0.1158956859312148
$ws.Cells.Item(14,10).Value = 0.08036858123741962
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 0.169654
$ws.Cells.Item(14,14).Value = 0.508962
$ws.Cells.Item(14,15).Value = 0.006094264463659866
$ws.Cells.Item(14,16).Value = 0.006534681579452628
$ws.Cells.Item(14,17).Value = 0.039688093317
$ws.Cells.Item(14,18).Value = 0.238128559902
$ws.Cells.Item(14,19).Value = 0.0007062989602620868
$ws.Cells.Item(14,20).Value = 0.0005251830873789081

# Row 15: MuSCs -> FAPs
$ws.Cells.Item(15,1).Value = "MuSCs"
$ws.Cells.Item(15,2).Value = "Slit1"
$ws.Cells.Item(15,3).Value = "Robo1"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 2
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 0.2339355
$ws.Cells.Item(15,8).Value = 0.467871
$ws.Cells.Item(15,9).Value = 0.1158956859312148
$ws.Cells.Item(15,10).Value = 0.08036858123741962
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 21.99231
$ws.Cells.Item(15,14).Value = 65.97693
$ws.Cells.Item(15,15).Value = 0.7900017288527916
$ws.Cells.Item(15,16).Value = 0.8470931604713817
$ws.Cells.Item(15,17).Value = 5.144782036005
$ws.Cells.Item(15,18).Value = 30.86869221603
$ws.Cells.Item(15,19).Value = 0.09155779225223983
$ws.Cells.Item(15,20).Value = 0.06807967548300678

# Row 16: MuSCs -> Inflammatory-Mac
$ws.Cells.Item(16,1).Value = "MuSCs"
$ws.Cells.Item(16,2).Value = "Slit1"
$ws.Cells.Item(16,3).Value = "Robo1"
$ws.Cells.Item(16,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(16,5).Value = 2
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 0.2339355
$ws.Cells.Item(16,8).Value = 0.467871
$ws.Cells.Item(16,9).Value = 0.1158956859312148
$ws.Cells.Item(16,10).Value = 0.08036858123741962
$ws.Cells.Item(16,11).Value = 1
$ws.Cells.Item(16,12).Value = 0.3333333333333333
$ws.Cells.Item(16,13).Value = 0.04769766666666667
$ws.Cells.Item(16,14).Value = 0.143093
$ws.Cells.Item(16,15).Value = 0.001713382501834088
$ws.Cells.Item(16,16).Value = 0.001837204332049573
$ws.Cells.Item(16,17).Value = 0.0111581775005
$ws.Cells.Item(16,18).Value = 0.06694906500299999
$ws.Cells.Item(16,19).Value = 0.0001985736403126025
$ws.Cells.Item(16,20).Value = 0.0001476535056100653

# Row 17: MuSCs -> MuSCs
$ws.Cells.Item(17,1).Value = "MuSCs"
$ws.Cells.Item(17,2).Value = "Slit1"
$ws.Cells.Item(17,3).Value = "Robo1"
$ws.Cells.Item(17,4).Value = "MuSCs"
$ws.Cells.Item(17,5).Value = 2
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 0.2339355
$ws.Cells.Item(17,8).Value = 0.467871
$ws.Cells.Item(17,9).Value = 0.1158956859312148
$ws.Cells.Item(17,10).Value = 0.08036858123741962
$ws.Cells.Item(17,11).Value = 2
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 5.6286445
$ws.Cells.Item(17,14).Value = 11.257289
$ws.Cells.Item(17,15).Value = 0.2021906241817143
$ws.Cells.Item(17,16).Value = 0.1445349536171162
$ws.Cells.Item(17,17).Value = 1.31673976542975
$ws.Cells.Item(17,18).Value = 5.266959061719
$ws.Cells.Item(17,19).Value = 0.02343302107840024
$ws.Cells.Item(17,20).Value = 0.01161606916142388

# Row 18: Resolving-Mac -> ECs
$ws.Cells.Item(18,1).Value = "Resolving-Mac"
$ws.Cells.Item(18,2).Value = "Slit1"
$ws.Cells.Item(18,3).Value = "Robo1"
$ws.Cells.Item(18,4).Value = "ECs"
$ws.Cells.Item(18,5).Value = 3
$ws.Cells.Item(18,6).Value = 1
$ws.Cells.Item(18,7).Value = 0.9162119999999999
$ws.Cells.Item(18,8).Value = 2.748636
$ws.Cells.Item(18,9).Value = 0.4539072445114579
$ws.Cells.Item(18,10).Value = 0.4721471851388441
$ws.Cells.Item(18,11).Value = 3
$ws.Cells.Item(18,12).Value = 1
$ws.Cells.Item(18,13).Value = 0.169654
$ws.Cells.Item(18,14).Value = 0.508962
$ws.Cells.Item(18,15).Value = 0.006094264463659866
$ws.Cells.Item(18,16).Value = 0.006534681579452628
$ws.Cells.Item(18,17).Value = 0.155439030648
$ws.Cells.Item(18,18).Value = 1.398951275832
$ws.Cells.Item(18,19).Value = 0.002766230790023947
$ws.Cells.Item(18,20).Value = 0.003085331513517214

# Row 19: Resolving-Mac -> FAPs
$ws.Cells.Item(19,1).Value = "Resolving-Mac"
$ws.Cells.Item(19,2).Value = "Slit1"
$ws.Cells.Item(19,3).Value = "Robo1"
$ws.Cells.Item(19,4).Value = "FAPs"
$ws.Cells.Item(19,5).Value = 3
$ws.Cells.Item(19,6).Value = 1
$ws.Cells.Item(19,7).Value = 0.9162119999999999
$ws.Cells.Item(19,8).Value = 2.748636
$ws.Cells.Item(19,9).Value = 0.4539072445114579
$ws.Cells.Item(19,10).Value = 0.4721471851388441
$ws.Cells.Item(19,11).Value = 3
$ws.Cells.Item(19,12).Value = 1
$ws.Cells.Item(19,13).Value = 21.99231
$ws.Cells.Item(19,14).Value = 65.97693
$ws.Cells.Item(19,15).Value = 0.7900017288527916
$ws.Cells.Item(19,16).Value = 0.8470931604713817
$ws.Cells.Item(19,17).Value = 20.14961832972
$ws.Cells.Item(19,18).Value = 181.34656496748
$ws.Cells.Item(19,19).Value = 0.3585875079028585
$ws.Cells.Item(19,20).Value = 0.39995265126693

# Row 20: Resolving-Mac -> Inflammatory-Mac
$ws.Cells.Item(20,1).Value = "Resolving-Mac"
$ws.Cells.Item(20,2).Value = "Slit1"
$ws.Cells.Item(20,3).Value = "Robo1"
$ws.Cells.Item(20,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(20,5).Value = 3
$ws.Cells.Item(20,6).Value = 1
$ws.Cells.Item(20,7).Value = 0.9162119999999999
$ws.Cells.Item(20,8).Value = 2.748636
$ws.Cells.Item(20,9).Value = 0.4539072445114579
$ws.Cells.Item(20,10).Value = 0.4721471851388441
$ws.Cells.Item(20,11).Value = 1
$ws.Cells.Item(20,12).Value = 0.3333333333333333
$ws.Cells.Item(20,13).Value = 0.04769766666666667
$ws.Cells.Item(20,14).Value = 0.143093
$ws.Cells.Item(20,15).Value = 0.001713382501834088
$ws.Cells.Item(20,16).Value = 0.001837204332049573
$ws.Cells.Item(20,17).Value = 0.043701174572
$ws.Cells.Item(20,18).Value = 0.393310571148
$ws.Cells.Item(20,19).Value = 0.0007777167302016589
$ws.Cells.Item(20,20).Value = 0.0008674308539020961

# Row 21: Resolving-Mac -> MuSCs
$ws.Cells.Item(21,1).Value = "Resolving-Mac"
$ws.Cells.Item(21,2).Value = "Slit1"
$ws.Cells.Item(21,3).Value = "Robo1"
$ws.Cells.Item(21,4).Value = "MuSCs"
$ws.Cells.Item(21,5).Value = 3
$ws.Cells.Item(21,6).Value = 1
$ws.Cells.Item(21,7).Value = 0.9162119999999999
$ws.Cells.Item(21,8).Value = 2.748636
$ws.Cells.Item(21,9).Value = 0.4539072445114579
$ws.Cells.Item(21,10).Value = 0.4721471851388441
$ws.Cells.Item(21,11).Value = 2
$ws.Cells.Item(21,12).Value = 1
$ws.Cells.Item(21,13).Value = 5.6286445
$ws.Cells.Item(21,14).Value = 11.257289
$ws.Cells.Item(21,15).Value = 0.2021906241817143
$ws.Cells.Item(21,16).Value = 0.1445349536171162
$ws.Cells.Item(21,17).Value = 5.157031634633999
$ws.Cells.Item(21,18).Value = 30.942189807804
$ws.Cells.Item(21,19).Value = 0.09177578908837368
$ws.Cells.Item(21,20).Value = 0.0682417715044948

